$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report date range) ---
$ws.Range("A8").Value = "Volume 30   Number  51"
$ws.Range("C9").Value = "Report Covering the Week  12/18/2023  Through  12/24/2023"

# --- Cells that change from numeric values to text placeholders ("0" / "***.*") ---
# C14 (text "0", style 14) and E14 (text "***.*", style 14) are untouched by this
# edit elsewhere, so copying their value+format over is a safe way to turn the
# destination cells into the same text-styled placeholder cells.
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("C14").Copy($ws.Range("F15"))
$ws.Range("C14").Copy($ws.Range("D26"))
$ws.Range("C14").Copy($ws.Range("F26"))
$ws.Range("C14").Copy($ws.Range("F30"))
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("E14").Copy($ws.Range("E26"))

# --- Plain numeric value updates (weekly crime-statistics refresh) ---
$ws.Range("N14").Value = -77.272727272727
$ws.Range("H15").Value = -100
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 200
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = 28.571428571428
$ws.Range("I16").Value = 194
$ws.Range("J16").Value = 175
$ws.Range("K16").Value = 10.857142857142
$ws.Range("L16").Value = 42.647058823529
$ws.Range("M16").Value = -46.703296703296
$ws.Range("N16").Value = -84.712371946414
$ws.Range("C17").Value = 4
$ws.Range("E17").Value = 300
$ws.Range("F17").Value = 31
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = 244.444444444444
$ws.Range("I17").Value = 310
$ws.Range("J17").Value = 266
$ws.Range("K17").Value = 16.541353383458
$ws.Range("L17").Value = 44.186046511627
$ws.Range("M17").Value = 51.219512195122
$ws.Range("N17").Value = -45.993031358885
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -57.142857142857
$ws.Range("F18").Value = 18
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = -5.263157894736
$ws.Range("I18").Value = 274
$ws.Range("J18").Value = 294
$ws.Range("K18").Value = -6.802721088435
$ws.Range("L18").Value = 18.614718614718
$ws.Range("M18").Value = -40.820734341252
$ws.Range("N18").Value = -79.131759329779
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 40
$ws.Range("F19").Value = 71
$ws.Range("G19").Value = 39
$ws.Range("H19").Value = 82.051282051282
$ws.Range("I19").Value = 775
$ws.Range("J19").Value = 652
$ws.Range("K19").Value = 18.865030674846
$ws.Range("L19").Value = 50.193798449612
$ws.Range("M19").Value = 54.690618762475
$ws.Range("N19").Value = 41.423357664233
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 200
$ws.Range("F20").Value = 14
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 186
$ws.Range("J20").Value = 181
$ws.Range("K20").Value = 2.762430939226
$ws.Range("L20").Value = 22.368421052631
$ws.Range("M20").Value = 15.527950310559
$ws.Range("N20").Value = -79.848320693391
$ws.Range("C21").Value = 27
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = 35
$ws.Range("F21").Value = 152
$ws.Range("G21").Value = 96
$ws.Range("H21").Value = 58.333333333333
$ws.Range("I21").Value = 1756
$ws.Range("J21").Value = 1584
$ws.Range("K21").Value = 10.858585858585
$ws.Range("L21").Value = 37.725490196078
$ws.Range("M21").Value = 3.051643192488
$ws.Range("N21").Value = -62.502669229126
$ws.Range("C22").Value = 2
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 26
$ws.Range("K22").Value = -33.333333333333
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -38.095238095238
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 200
$ws.Range("G23").Value = 10
$ws.Range("H23").Value = 90
$ws.Range("I23").Value = 204
$ws.Range("J23").Value = 156
$ws.Range("K23").Value = 30.769230769230
$ws.Range("L23").Value = 20.710059171597
$ws.Range("M23").Value = 55.725190839694
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = -11.111111111111
$ws.Range("F24").Value = 88
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 1026
$ws.Range("J24").Value = 1174
$ws.Range("K24").Value = -12.606473594548
$ws.Range("L24").Value = 0.984251968503
$ws.Range("M24").Value = -17.191283292978
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 166.666666666667
$ws.Range("F25").Value = 32
$ws.Range("G25").Value = 27
$ws.Range("H25").Value = 18.518518518518
$ws.Range("I25").Value = 480
$ws.Range("J25").Value = 487
$ws.Range("K25").Value = -1.437371663244
$ws.Range("L25").Value = 29.032258064516
$ws.Range("M25").Value = -0.826446280991
$ws.Range("H26").Value = -100
$ws.Range("C27").Value = 2
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 200
$ws.Range("I27").Value = 53
$ws.Range("K27").Value = 15.217391304347
$ws.Range("L27").Value = -20.895522388059
